$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted above the existing row 418
# (Especial quality, Frutilla, Terminal Hortofrutícola Agro Chillán),
# pushing all subsequent records (old rows 418-509) down by one row.
$ws.Rows.Item(418).Insert()

# Populate the newly inserted row 418 with the new observation.
$ws.Range("A418").Value2 = 7
$ws.Range("B418").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C418").Value2 = "Ñuble"
$ws.Range("D418").Value2 = 45015
$ws.Range("E418").Value2 = 16
$ws.Range("F418").Value2 = "Fruta"
$ws.Range("G418").Value2 = 100101
$ws.Range("H418").Value2 = "Berries"
$ws.Range("I418").Value2 = 100112025
$ws.Range("J418").Value2 = "Frutilla"
$ws.Range("K418").Value2 = "Sin especificar"
$ws.Range("L418").Value2 = "Especial"
$ws.Range("M418").Value2 = 40
$ws.Range("N418").Value2 = 7000
$ws.Range("O418").Value2 = 7000
$ws.Range("P418").Value2 = 7000
$ws.Range("Q418").Value2 = "`$/caja 7 kilos"
$ws.Range("R418").Value2 = "Provincia de Diguillín"
$ws.Range("S418").Value2 = 1000
$ws.Range("T418").Value2 = 7
